# Apply the two changes described by the commit:
#  1. Fixed "date updated automatically" placeholder text moves from
#     24/8/2017 -> 26/8/2017 on the slide master and every slide layout.
#  2. Typo fix on slide 2: "PRECENTER" -> "PRESENTER".

$p = $ppt.ActivePresentation

$oldDate = "24/8/2017"
$newDate = "26/8/2017"

# --- 1. Slide master date placeholder -------------------------------------
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- 2. Every slide layout's date placeholder ------------------------------
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- 3. Fix "PRECENTER" -> "PRESENTER" on slide 2 --------------------------
$s2 = $p.Slides.Item(2)
$grp = $s2.Shapes.Item(1)
for ($i = 1; $i -le $grp.GroupItems.Count; $i++) {
    $sh = $grp.GroupItems.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq "PRECENTER") {
            $sh.TextFrame.TextRange.Text = "PRESENTER"
        }
    }
}
